$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Clear out the old content of rows 14-18 (columns A-H) so we can rebuild
# them to match the new layout (rows 14-17).
$ws.Range("A14:H18").ClearContents()

# Row 14 (PlayerTest testcase 3): bet <= 0 (B) and players* (E)
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = "x"
$ws.Range("E14").Value = "x"

# Row 15 (PlayerTest testcase 4): bet > chips (C) and players = null (D)
$ws.Range("A15").Value = 4
$ws.Range("C15").Value = "x"
$ws.Range("D15").Value = "x"

# Row 16 (PlayerTest testcase 5)
$ws.Range("A16").Value = 5
$ws.Range("C16").Value = "x"
$ws.Range("E16").Value = "x"
$ws.Range("G16").Value = "x"

# Row 17: comment row (moved up from old row 18)
$ws.Range("A17").Value = "// Vi kommer behöva förklara att inga parametrar är en valid ekvivalensklass också"

# Update the selection to reflect the new active cell
$ws.Range("E14").Select()
